$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "FER60" sample label carried a stray trailing (non-breaking) space in
# the shared-string table; clean it up now that the annulus-count error
# estimates are finalized. The row/column data (FER60 -> recount 3) is
# unchanged, only the label text itself.
$ws.Range("A11").Value = "FER60"

# Restore the author's final cursor position in the sheet.
$ws.Range("D17").Select()
